$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user_detail rows appended to the master data table (rows 22-30).
$rows = @(
    @(110021, 7316931025, "Magdalena Weber",   "magdalena.weber@xyz.com",   932122450),
    @(110022, 9137847236, "Adrienne Hoffman",  "adrienne.hoffman@xyz.com",  848488000),
    @(110023, 8428758532, "Adrienne Mcgee",    "adrienne.mcgee@xyz.com",    894773246),
    @(110024, 9804209494, "Amare Coleman",     "amare.coleman@xyz.com",     956554588),
    @(110025, 7105248214, "Dawson Ibarra",     "dawson.ibarra@xyz.com",     765455583),
    @(110026, 9316557128, "Elvis Mcmillan",    "elvis.mcmillan@xyz.com",    884282274),
    @(110027, 8103486949, "Steve George",      "steve.george@xyz.com",      971073663),
    @(110028, 9601932866, "Colton Elliott",    "colton.elliott@xyz.com",    809908673),
    @(110029, 9317596765, "Carolyn Rodriguez", "carolyn.rodriguez@xyz.com", 818876429)
)

$startRow = 22
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = "ACT"
    $ws.Cells.Item($r, 7).Value = "eng"
    $ws.Cells.Item($r, 8).Value = "PWD"
    $ws.Cells.Item($r, 9).Value = $true
    $ws.Cells.Item($r, 10).Value = "superadmin"
    $ws.Cells.Item($r, 11).Value = "now()"
    $ws.Cells.Item($r, 12).Value = "now()"
}

# Match the look of the existing table: column D (email) and column I
# (is_active) carry formatting beyond the default, same as rows 2:21.
$ws.Range("D2").Copy()
$ws.Range("D22:D30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I2").Copy()
$ws.Range("I22:I30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the viewport / selection the way the author left it after adding data.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("A22:A30").Select() | Out-Null
